$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update building/shop names in column A (rows 2-10) - translate codes to Russian names
$ws.Range("A2").Value = "ЮАБЗ"
$ws.Range("A3").Value = "САБЗ"
$ws.Range("A4").Value = "ЗУ"
$ws.Range("A5").Value = "ЦЗЛ"
$ws.Range("A6").Value = "Инженерный корпус"
$ws.Range("A7").Value = "24 цех"
$ws.Range("A10").Value = "28 цех"
$ws.Range("A9").Value = "202 цех"
$ws.Range("A8").Value = "24 цех АБЗ"

# Update selection to A6
$ws.Range("A6").Select()
